$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously had 6 data rows (rows 2-7): every combination of
# Sending cluster (ECs/FAPs/MuSCs) x Target cluster (ECs/MuSCs). The new
# TPM-based run only keeps the Target cluster = MuSCs rows, so drop the
# three Target cluster = ECs rows. Deleting "row 2" three times in a row
# removes original rows 2, 4 and 6 because each delete shifts the rows
# below it up by one.
$ws.Rows(2).Delete()
$ws.Rows(3).Delete()
$ws.Rows(4).Delete()

# Refresh the surviving rows with the recomputed TPM-derived values.
# Row 2: Sending cluster = ECs, Target cluster = MuSCs
$ws.Range("G2").Value = 0.5347833333333334
$ws.Range("H2").Value = 1.60435
$ws.Range("I2").Value = 0.196822066153855
$ws.Range("J2").Value = 0.196822066153855
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.0257722784
$ws.Range("R2").Value = 0.2319505056000001
$ws.Range("S2").Value = 0.196822066153855
$ws.Range("T2").Value = 0.196822066153855

# Row 3: Sending cluster = FAPs, Target cluster = MuSCs
$ws.Range("I3").Value = 0.1891972429821067
$ws.Range("J3").Value = 0.1891972429821067
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("S3").Value = 0.1891972429821067
$ws.Range("T3").Value = 0.1891972429821067

# Row 4: Sending cluster = MuSCs, Target cluster = MuSCs
$ws.Range("I4").Value = 0.6139806908640383
$ws.Range("J4").Value = 0.6139806908640382
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("S4").Value = 0.6139806908640383
$ws.Range("T4").Value = 0.6139806908640382
